$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting the old row 47 (and below) down to 48.
$ws.Rows.Item(47).Insert()

# The new row 47 receives a copy of what used to be in row 46 (before its values
# are updated below), i.e. the 2021-08-02 / 400 / 14000 / 15000 / 14575 / 1121 record.
$ws.Range("A47").Value = 6
$ws.Range("B47").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 44410
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = 100114007
$ws.Range("G47").Value = "Jengibre"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 400
$ws.Range("K47").Value = 14000
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = 14575
$ws.Range("N47").Value = "$/caja 13 kilos"
$ws.Range("O47").Value = "Perú"
$ws.Range("P47").Value = 1121
$ws.Range("Q47").Value = 13
$ws.Range("R47").Value = "Hortaliza"

# Row 46 itself is updated to the new weekly record (2021-09-22).
$ws.Range("D46").Value = 44461
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 12000
$ws.Range("L46").Value = 13000
$ws.Range("M46").Value = 12400
$ws.Range("P46").Value = 954
